$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.04'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.47%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.20%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.131'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.31%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05604'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.39%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.474'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.02%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8191'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.08%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8326'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.96%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1331'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.14%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06949'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.05%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02895'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.69%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09385'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.15%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001511'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.17%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006012'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.04%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006243'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.39%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '3.71%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.024'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.27%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '8.29%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.12%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03081'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-4.60%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.21%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.758'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.31%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04601'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-2.27%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-1.67%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004491'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-2.59%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009607'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.96%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001400'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '0.70%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.43%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1375'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '30.69%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006224'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.38%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002602'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.35%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009030'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '18.17%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005352'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.73%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.03%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1400'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '4.89%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002458'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '15.78%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.03%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.03%'
